$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 1000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -1936
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1468
$ws.Range("H28").Value = 689.6
$ws.Range("I28").Value = 512
$ws.Range("K28").Value = 512
$ws.Range("M28").Value = -27
$ws.Range("H43").Value = 6831.6665
$ws.Range("I43").Value = 3747
$ws.Range("J43").Value = 13001
$ws.Range("K43").Value = 3747
$ws.Range("L43").Value = 13001
$ws.Range("M43").Value = -3678
$ws.Range("N43").Value = -13139
$ws.Range("H55").Value = 633.3333
$ws.Range("I55").Value = 1000
$ws.Range("K55").Value = 1000
$ws.Range("M55").Value = -786
$ws.Range("H58").Value = 321
$ws.Range("J58").Value = 397.5
$ws.Range("L58").Value = 1192.5
$ws.Range("N58").Value = -1492.5
$ws.Range("H127").Value = 4659.6
$ws.Range("J127").Value = 2500
$ws.Range("L127").Value = 7500
$ws.Range("N127").Value = -17420

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 107.666664
$ws.Range("I5").Value = 99
$ws.Range("K5").Value = 99
$ws.Range("M5").Value = 13
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2288
$ws.Range("H96").Value = 28061.334
$ws.Range("J96").Value = 28061.334
$ws.Range("L96").Value = 28061.334
$ws.Range("N96").Value = -33553.334
$ws.Range("H97").Value = 545.6
$ws.Range("J97").Value = 299
$ws.Range("L97").Value = 299
$ws.Range("N97").Value = -1291
$ws.Range("H122").Value = 1442
$ws.Range("I122").Value = 1256
$ws.Range("K122").Value = 3768
$ws.Range("M122").Value = -1318
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 107.666664
$ws.Range("I4").Value = 99
$ws.Range("K4").Value = 99
$ws.Range("M4").Value = 16
$ws.Range("H20").Value = 8682
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H94").Value = 503.14285
$ws.Range("I94").Value = 571
$ws.Range("K94").Value = 571
$ws.Range("M94").Value = -120
$ws.Range("H99").Value = 1950
$ws.Range("I99").Value = 1950
$ws.Range("K99").Value = 1950
$ws.Range("M99").Value = -452
$ws.Range("H105").Value = 2674.3635
$ws.Range("I105").Value = 2614.3
$ws.Range("J105").Value = 3275
$ws.Range("K105").Value = 2614.3
$ws.Range("L105").Value = 3275
$ws.Range("M105").Value = -867.3000000000002
$ws.Range("N105").Value = -6769

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.1
$ws.Range("I7").Value = 69.85714
$ws.Range("J7").Value = 24
$ws.Range("K7").Value = 69.85714
$ws.Range("L7").Value = 24
$ws.Range("M7").Value = 43.14286
$ws.Range("N7").Value = -250
$ws.Range("H22").Value = 839.8
$ws.Range("I22").Value = 233.33333
$ws.Range("K22").Value = 233.33333
$ws.Range("M22").Value = 116.66667
$ws.Range("H39").Value = 6750
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("N39").Value = -10782
$ws.Range("H45").Value = 12329.333
$ws.Range("I45").Value = 995
$ws.Range("J45").Value = 34998
$ws.Range("K45").Value = 995
$ws.Range("L45").Value = 34998
$ws.Range("M45").Value = -402
$ws.Range("N45").Value = -36184
$ws.Range("H49").Value = 6750
$ws.Range("J49").Value = 10000
$ws.Range("L49").Value = 10000
$ws.Range("N49").Value = -10364
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1024.9333
$ws.Range("I2").Value = 375
$ws.Range("J2").Value = 5249.5
$ws.Range("K2").Value = 2250
$ws.Range("L2").Value = 31497
$ws.Range("M2").Value = -2137
$ws.Range("N2").Value = -31723
$ws.Range("H13").Value = 3350
$ws.Range("I13").Value = 3133.3333
$ws.Range("K13").Value = 9399.999899999999
$ws.Range("M13").Value = -9231.999899999999
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("K63").Value = 9000
$ws.Range("M63").Value = -8251
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("K66").Value = 27000
$ws.Range("M66").Value = -23256
$ws.Range("H103").Value = 4466.1665
$ws.Range("I103").Value = 2000
$ws.Range("K103").Value = 6000
$ws.Range("M103").Value = -5121
$ws.Range("H109").Value = 17887.5
$ws.Range("I109").Value = 17887.5
$ws.Range("K109").Value = 53662.5
$ws.Range("M109").Value = -52622.5
$ws.Range("H139").Value = 2758
$ws.Range("I139").Value = 2758
$ws.Range("K139").Value = 8274
$ws.Range("M139").Value = -3134
$ws.Range("H140").Value = 3082
$ws.Range("I140").Value = 3082
$ws.Range("K140").Value = 9246
$ws.Range("M140").Value = -4066

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2373.1667
$ws.Range("J3").Value = 2059.75
$ws.Range("L3").Value = 2059.75
$ws.Range("N3").Value = -2291.75
$ws.Range("H15").Value = 23999
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 23999
$ws.Range("H84").Value = 23999
$ws.Range("H122").Value = 2000
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2119.8
$ws.Range("I7").Value = 1649.75
$ws.Range("K7").Value = 1649.75
$ws.Range("M7").Value = -1537.75
$ws.Range("H22").Value = 2172.2
$ws.Range("I22").Value = 2629.8
$ws.Range("J22").Value = 1714.6
$ws.Range("K22").Value = 2629.8
$ws.Range("L22").Value = 1714.6
$ws.Range("M22").Value = -2334.8
$ws.Range("N22").Value = -2304.6
$ws.Range("H27").Value = 2172.2
$ws.Range("I27").Value = 2629.8
$ws.Range("J27").Value = 1714.6
$ws.Range("K27").Value = 2629.8
$ws.Range("L27").Value = 1714.6
$ws.Range("M27").Value = -2522.8
$ws.Range("N27").Value = -1928.6
$ws.Range("H40").Value = 999
$ws.Range("I40").Value = 999
$ws.Range("K40").Value = 999
$ws.Range("M40").Value = -863
$ws.Range("H82").Value = 1725.9
$ws.Range("I82").Value = 899.2857
$ws.Range("J82").Value = 3654.6667
$ws.Range("K82").Value = 899.2857
$ws.Range("L82").Value = 3654.6667
$ws.Range("M82").Value = -538.2857
$ws.Range("N82").Value = -4376.6667
$ws.Range("H85").Value = 1725.9
$ws.Range("I85").Value = 899.2857
$ws.Range("J85").Value = 3654.6667
$ws.Range("K85").Value = 899.2857
$ws.Range("L85").Value = 3654.6667
$ws.Range("M85").Value = 348.7143
$ws.Range("N85").Value = -6150.6667
$ws.Range("H122").Value = 2517
$ws.Range("I122").Value = 1197
$ws.Range("J122").Value = 3045
$ws.Range("K122").Value = 3591
$ws.Range("L122").Value = 9135
$ws.Range("M122").Value = -1141
$ws.Range("N122").Value = -14035
$ws.Range("H126").Value = 2119.8
$ws.Range("I126").Value = 1649.75
$ws.Range("K126").Value = 4949.25
$ws.Range("M126").Value = -2479.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 996.5
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 996.5
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 15000
$ws.Range("N84").Value = -25608
$ws.Range("H95").Value = 39991
$ws.Range("J95").Value = 39991
$ws.Range("L95").Value = 39991
$ws.Range("N95").Value = -45483
$ws.Range("H122").Value = 5113.75
$ws.Range("I122").Value = 4966.6665
$ws.Range("K122").Value = 14899.9995
$ws.Range("M122").Value = -12449.9995
